# Add a new row (row 94) of index data to "Sheet 1", following the same
# pattern as the existing rows: a date in column A (formatted like the
# dates above it), plain numeric index values in B/C, and index values
# stored as text in D:G (mirrors how the source data file was produced).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

# Column A: new date. Copy the cell above so the existing date style
# (xf with numFmtId 14) is reused instead of a brand-new number format
# being registered, then overwrite the copied value.
$ws.Range("A93").Copy($ws.Range("A94"))
$ws.Range("A94").Value = 45536

# Columns B and C: plain numeric index values.
$ws.Range("B94").Value = 105.287252600406
$ws.Range("C94").Value = 119.922839673779

# Columns D-G: index values stored as literal text (not numbers), matching
# the source file. Assigning a numeric-looking string straight to .Value
# would be auto-coerced to a number (like typing it into Excel), so each
# value is first built as a text formula result in a scratch cell (which
# keeps it a genuine string) and then copied into place, which carries the
# string - and only the string - into the destination cell/style.
function Set-TextValue($cell, $text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy($cell)
    $scratch.ClearContents()
}

Set-TextValue $ws.Range("D94") "110.4"
Set-TextValue $ws.Range("E94") "111.5"
Set-TextValue $ws.Range("F94") " 88.3"
Set-TextValue $ws.Range("G94") "168.5"
